# Update report period labels (row 9, columns G and H)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G9").Value = "1402-04-06 (11)"
$ws.Range("H9").Value = "1402-04-06 (3)"

# Update latest-period (column H) figures with refreshed data
$ws.Range("H14").Value = -27471
$ws.Range("H17").Value = 35823
$ws.Range("H19").Value = 77462
$ws.Range("H20").Value = 91124
$ws.Range("H21").Value = -2001
$ws.Range("H22").Value = 89123
$ws.Range("H24").Value = 89123
